# Apply the "consolidated report" updates to the Absent (H) column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that were marked Absent (H) = 1
$ws.Range("H3").Value  = 1
$ws.Range("H9").Value  = 1
$ws.Range("H11").Value = 1

# Rows whose Absent (H) cell was an empty placeholder, now filled with 0
$ws.Range("H5").Value  = 0
$ws.Range("H10").Value = 0
$ws.Range("H12").Value = 0
